# Bug fix in Eduati data files:
# Sheet1 ("SNUC2B_noCTRL_meas") had 43 extra leftover rows (45:87) that only
# carried a stray incrementing index in column A - remove them so the sheet's
# used range matches the real 44-row dataset (same shape as Sheet2 / Sheet3).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Drop the stray rows 45:87 on Sheet1 - they only had a leftover value in
# column A, no real data - collapsing the used range back down to A1:N44.
[void]$ws1.Rows("45:87").Delete()

# Sheet1 becomes the active/selected sheet (it had been Sheet3), scrolled to
# the bottom of the data and with F61 (now outside the data, below row 44)
# selected - mirrors the saved view state after the cleanup.
[void]$ws1.Activate()
[void]$ws1.Range("F61").Select()

# Sheet3 is no longer the active tab - its selection goes back to the default
# top-left data block.
[void]$ws3.Range("A2:N44").Select()

# Leave Sheet1 as the active sheet/tab.
[void]$ws1.Activate()
